$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.085.27'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '2.090.73'
$ws.Range('E3').Value = '  +2.68%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.16'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.614'
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '61.09'
$ws.Range('E7').Value = '  +1.29%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.382'
$ws.Range('E9').Value = '  -0.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0841'
$ws.Range('E10').Value = '  +2.71%  '
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').Value = '2.402.71'
$ws.Range('E12').Value = '  +2.79%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.60'
$ws.Range('E13').Value = '  +0.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.09'
$ws.Range('E14').Value = '  +3.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.51'
$ws.Range('E15').Value = '  +6.61%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.775'
$ws.Range('E16').Value = '  +1.62%  '
$ws.Range('D17').Value = '2.089.05'
$ws.Range('E17').Value = '  +2.59%  '
$ws.Range('D18').Value = '38.027.05'
$ws.Range('E18').Value = '  +0.39%  '
$ws.Range('E19').Value = '  +1.90%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.15'
$ws.Range('E20').Value = '  +0.29%  '
$ws.Range('E21').Value = '  +1.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '224.23'
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('E23').Value = '  +0.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.43'
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('E25').Value = '  +3.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '169.48'
$ws.Range('E26').Value = '  +1.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.40'
$ws.Range('E27').Value = '  +0.91%  '
$ws.Range('E28').Value = '  +2.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.00'
$ws.Range('E29').Value = '  +0.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.32'
$ws.Range('E30').Value = '  +3.39%  '
$ws.Range('E31').Value = '  -0.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.35'
$ws.Range('E32').Value = '  +9.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.44'
$ws.Range('E33').Value = '  +0.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.66'
$ws.Range('E34').Value = '  +3.13%  '
$ws.Range('E35').Value = '  +0.18%  '
$ws.Range('B36').Value = 'THORChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.49'
$ws.Range('E36').Value = '  +0.56%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.41'
$ws.Range('E37').Value = '  +6.00%  '
$ws.Range('E38').Value = '  +8.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.91'
$ws.Range('E40').Value = '  +4.48%  '
$ws.Range('D41').Value = '1.546.66'
$ws.Range('E41').Value = '  +1.39%  '
$ws.Range('E42').Value = '  +4.25%  '
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('E44').Value = '  -0.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0903'
$ws.Range('E45').Value = '  -1.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.16'
$ws.Range('E46').Value = '  +3.82%  '
$ws.Range('E47').Value = '  +0.75%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.00'
$ws.Range('E49').Value = '  +0.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.22'
$ws.Range('E50').Value = '  +1.48%  '
$ws.Range('D51').Value = '2.289.73'
$ws.Range('E51').Value = '  +2.84%  '
